# Apply the latest scraped cryptocurrency prices / 1h volume changes
# (and re-sort the Filecoin / ImmutableX pair) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update D, E
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '22.290.14'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -5.00%  '

# Row 3: update D, E
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.560.31'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.19%  '

# Row 5: update D, E
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '1.002'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.12%  '

# Row 6: update D, E
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '288.64'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.10%  '

# Row 7: update D, E
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3726'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.72%  '

# Row 8: update E
$ws.Range('E8').Value = '  -2.45%  '

# Row 9: update D, E
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3398'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.11%  '

# Row 10: update E
$ws.Range('E10').Value = '  -4.47%  '

# Row 11: update D, E
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07630'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.46%  '

# Row 12: update D, E
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.01%  '

# Row 13: update D, E
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.31'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.56%  '

# Row 14: update D, E
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.020'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.43%  '

# Row 15: update D, E
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.908'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -4.80%  '

# Row 16: update D, E
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.561.49'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.26%  '

# Row 17: update D, E
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001125'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -7.15%  '

# Row 18: update D, E
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '89.69'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.87%  '

# Row 19: update D, E
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06718'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.96%  '

# Row 20: update D, E
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.02%  '

# Row 21: update D, E
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.232'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.11%  '

# Row 22: update D, E
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '16.54'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.07%  '

# Row 23: update D, E
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.5278'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -7.91%  '

# Row 24: update D, E
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.96'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -3.86%  '

# Row 25: update D, E
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '22.338.94'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.85%  '

# Row 26: update D, E
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.398'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.75%  '

# Row 27: update D, E
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.801'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -6.46%  '

# Row 28: update E
$ws.Range('E28').Value = '  -4.06%  '

# Row 29: update D, E
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '146.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.94%  '

# Row 30: update D, E
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.981'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.05%  '

# Row 31: update D, E
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '125.10'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.04%  '

# Row 32: update D
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.732.80'
$ws.Range('D32').Style = 'Normal'

# Row 33: update B, C, D, E
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.008'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.96%  '

# Row 34: update B, C, D, E
$ws.Range('B34').Value = 'Filecoin'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.157'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -9.96%  '

# Row 35: update D, E
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.011'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.11%  '

# Row 36: update D, E
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -10.43%  '

# Row 37: update D, E
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.08439'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.00%  '

# Row 38: update D, E
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02543'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.42%  '

# Row 39: update D, E
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2307'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.62%  '

# Row 40: update D, E
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.480'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.26%  '

# Row 41: update D, E
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.308'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.04%  '

# Row 42: update D, E
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.06379'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.06%  '

# Row 43: update D, E
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.66'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -9.54%  '

# Row 44: update D, E
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6339'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.91%  '

# Row 45: update D, E
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.000'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.00%  '

# Row 46: update E
$ws.Range('E46').Value = '  -9.70%  '

# Row 47: update D, E
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5959'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.74%  '

# Row 48: update E
$ws.Range('E48').Value = '  -4.38%  '

# Row 49: update D, E
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.089'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.88%  '

# Row 50: update D, E
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.263'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.93%  '

# Row 51: update D, E
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '124.30'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.51%  '
